# Scheduled-runner update: refresh cached market-board derived figures
# (currentAveragePrice*, Leve price/profit columns) across the per-job
# leve-profit sheets. Source data values only; no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1927242.9
$ws.Range("J43").Value = 4721.3335
$ws.Range("L43").Value = 4721.3335
$ws.Range("N43").Value = -4859.3335
$ws.Range("H55").Value = 634.1
$ws.Range("I55").Value = 1100
$ws.Range("K55").Value = 1100
$ws.Range("M55").Value = -886
$ws.Range("H116").Value = 57878804
$ws.Range("J116").Value = 100008800
$ws.Range("L116").Value = 100008800
$ws.Range("N116").Value = -100015684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20632
$ws.Range("I61").Value = 26948
$ws.Range("K61").Value = 26948
$ws.Range("M61").Value = -26736
$ws.Range("H97").Value = 677.7778
$ws.Range("I97").Value = 571.4286
$ws.Range("J97").Value = 1050
$ws.Range("K97").Value = 571.4286
$ws.Range("L97").Value = 1050
$ws.Range("M97").Value = -75.42859999999996
$ws.Range("N97").Value = -2042
$ws.Range("H132").Value = 12511.407
$ws.Range("I132").Value = 15158.275
$ws.Range("K132").Value = 45474.825
$ws.Range("M132").Value = -42944.825
$ws.Range("H136").Value = 20632
$ws.Range("I136").Value = 26948
$ws.Range("K136").Value = 80844
$ws.Range("M136").Value = -78294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 237.25
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = -77
$ws.Range("H64").Value = 2315329.2
$ws.Range("I64").Value = 8333771
$ws.Range("K64").Value = 8333771
$ws.Range("M64").Value = -8333546
$ws.Range("H67").Value = 2315329.2
$ws.Range("I67").Value = 8333771
$ws.Range("K67").Value = 8333771
$ws.Range("M67").Value = -8332991
$ws.Range("H99").Value = 5209579
$ws.Range("I99").Value = 10417171
$ws.Range("J99").Value = 1987.5
$ws.Range("K99").Value = 10417171
$ws.Range("L99").Value = 1987.5
$ws.Range("M99").Value = -10415673
$ws.Range("N99").Value = -4983.5
$ws.Range("H105").Value = 62501828
$ws.Range("I105").Value = 100001620
$ws.Range("K105").Value = 100001620
$ws.Range("M105").Value = -99999873
$ws.Range("H107").Value = 2558
$ws.Range("I107").Value = 3574.8
$ws.Range("J107").Value = 863.3333
$ws.Range("K107").Value = 3574.8
$ws.Range("L107").Value = 863.3333
$ws.Range("M107").Value = -1654.8
$ws.Range("N107").Value = -4703.3333
$ws.Range("H134").Value = 1858.4445
$ws.Range("I134").Value = 1715.875
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 5147.625
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -2612.625
$ws.Range("N134").Value = -14067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5914.1904
$ws.Range("I31").Value = 1994
$ws.Range("J31").Value = 6836.5884
$ws.Range("K31").Value = 1994
$ws.Range("L31").Value = 6836.5884
$ws.Range("M31").Value = -1699
$ws.Range("N31").Value = -7426.5884
$ws.Range("H34").Value = 5914.1904
$ws.Range("I34").Value = 1994
$ws.Range("J34").Value = 6836.5884
$ws.Range("K34").Value = 1994
$ws.Range("L34").Value = 6836.5884
$ws.Range("M34").Value = -1792
$ws.Range("N34").Value = -7240.5884
$ws.Range("H58").Value = 336815.1
$ws.Range("I58").Value = 770936.9
$ws.Range("J58").Value = 4839.647
$ws.Range("K58").Value = 770936.9
$ws.Range("L58").Value = 4839.647
$ws.Range("M58").Value = -770733.9
$ws.Range("N58").Value = -5245.647
$ws.Range("H132").Value = 15884989
$ws.Range("I132").Value = 20848344
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 62545032
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -62542502
$ws.Range("N132").Value = -11810
$ws.Range("H136").Value = 336815.1
$ws.Range("I136").Value = 770936.9
$ws.Range("J136").Value = 4839.647
$ws.Range("K136").Value = 2312810.7
$ws.Range("L136").Value = 14518.941
$ws.Range("M136").Value = -2310260.7
$ws.Range("N136").Value = -19618.941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 297.33334
$ws.Range("I14").Value = 297.33334
$ws.Range("K14").Value = 892.0000200000001
$ws.Range("M14").Value = -719.0000200000001
$ws.Range("H97").Value = 350.25
$ws.Range("I97").Value = 310.66666
$ws.Range("J97").Value = 374
$ws.Range("K97").Value = 931.9999799999999
$ws.Range("L97").Value = 1122
$ws.Range("M97").Value = -435.9999799999999
$ws.Range("N97").Value = -2114
$ws.Range("H122").Value = 1266.75
$ws.Range("I122").Value = 779.6667
$ws.Range("J122").Value = 1559
$ws.Range("K122").Value = 7017.0003
$ws.Range("L122").Value = 14031
$ws.Range("M122").Value = -4567.0003
$ws.Range("N122").Value = -18931
$ws.Range("H134").Value = 8715.526
$ws.Range("I134").Value = 2376.6667
$ws.Range("K134").Value = 7130.000100000001
$ws.Range("M134").Value = -2060.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2385932.5
$ws.Range("J70").Value = 4798.4287
$ws.Range("L70").Value = 4798.4287
$ws.Range("N70").Value = -5338.4287
$ws.Range("H73").Value = 2385932.5
$ws.Range("J73").Value = 4798.4287
$ws.Range("L73").Value = 4798.4287
$ws.Range("N73").Value = -6670.4287
$ws.Range("H97").Value = 562.0357
$ws.Range("I97").Value = 536.0952
$ws.Range("K97").Value = 536.0952
$ws.Range("M97").Value = -40.09519999999998
$ws.Range("H122").Value = 316972.16
$ws.Range("I122").Value = 395011.78
$ws.Range("J122").Value = 4813.7144
$ws.Range("K122").Value = 1185035.34
$ws.Range("L122").Value = 14441.1432
$ws.Range("M122").Value = -1182585.34
$ws.Range("N122").Value = -19341.1432
$ws.Range("H128").Value = 100780
$ws.Range("J128").Value = 100780
$ws.Range("L128").Value = 100780
$ws.Range("N128").Value = -110740
$ws.Range("H132").Value = 123672.586
$ws.Range("I132").Value = 185680.81
$ws.Range("K132").Value = 557042.4299999999
$ws.Range("M132").Value = -554512.4299999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5864.8423
$ws.Range("J46").Value = 5919.029
$ws.Range("L46").Value = 5919.029
$ws.Range("N46").Value = -6295.029
$ws.Range("H61").Value = 2425.8333
$ws.Range("I61").Value = 2259.6155
$ws.Range("K61").Value = 2259.6155
$ws.Range("M61").Value = -2057.6155
$ws.Range("H100").Value = 3950.6924
$ws.Range("I100").Value = 4101.375
$ws.Range("J100").Value = 3709.6
$ws.Range("K100").Value = 4101.375
$ws.Range("L100").Value = 3709.6
$ws.Range("M100").Value = -3560.375
$ws.Range("N100").Value = -4791.6
$ws.Range("H113").Value = 2425.8333
$ws.Range("I113").Value = 2259.6155
$ws.Range("K113").Value = 2259.6155
$ws.Range("M113").Value = -89.61549999999988
$ws.Range("H132").Value = 3302.1404
$ws.Range("I132").Value = 2506.3777
$ws.Range("J132").Value = 6286.25
$ws.Range("K132").Value = 7519.1331
$ws.Range("L132").Value = 18858.75
$ws.Range("M132").Value = -4989.1331
$ws.Range("N132").Value = -23918.75
$ws.Range("H136").Value = 3651.0527
$ws.Range("I136").Value = 2731.4
$ws.Range("J136").Value = 7099.75
$ws.Range("K136").Value = 8194.200000000001
$ws.Range("L136").Value = 21299.25
$ws.Range("M136").Value = -5644.200000000001
$ws.Range("N136").Value = -26399.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1386.875
$ws.Range("I113").Value = 1497.5
$ws.Range("J113").Value = 1276.25
$ws.Range("K113").Value = 4492.5
$ws.Range("L113").Value = 3828.75
$ws.Range("M113").Value = -2322.5
$ws.Range("N113").Value = -8168.75
$ws.Range("H126").Value = 2360.6086
$ws.Range("I126").Value = 1952.5834
$ws.Range("J126").Value = 2805.7273
$ws.Range("K126").Value = 5857.7502
$ws.Range("L126").Value = 8417.1819
$ws.Range("M126").Value = -3387.7502
$ws.Range("N126").Value = -13357.1819
$ws.Range("H132").Value = 42744080
$ws.Range("I132").Value = 9260342
$ws.Range("K132").Value = 27781026
$ws.Range("M132").Value = -27778496
$ws.Range("H135").Value = 34215
$ws.Range("J135").Value = 34215
$ws.Range("L135").Value = 34215
$ws.Range("N135").Value = -44355
